$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.10254410228201
$ws.Range("C2").Value = 9.890833230244272
$ws.Range("E2").Value = 14.36426845683459
$ws.Range("F2").Value = 41.44786610073778
$ws.Range("G2").Value = 43.35808356879687
$ws.Range("H2").Value = 17.74606095788246
$ws.Range("J2").Value = 8.544658145475138
$ws.Range("L2").Value = 12.30421110837015
$ws.Range("B3").Value = 19.61917466871718
$ws.Range("C3").Value = 9.596634703073315
$ws.Range("E3").Value = 14.38288654230258
$ws.Range("F3").Value = 41.46010470586681
$ws.Range("G3").Value = 43.32113652590436
$ws.Range("H3").Value = 17.80225474590215
$ws.Range("J3").Value = 8.571670516551224
$ws.Range("L3").Value = 12.2787717290679
$ws.Range("B4").Value = 19.32026502154306
$ws.Range("C4").Value = 9.409686018749261
$ws.Range("E4").Value = 14.39554085439247
$ws.Range("F4").Value = 41.47940035484071
$ws.Range("G4").Value = 43.31616473455267
$ws.Range("H4").Value = 17.84114128464582
$ws.Range("J4").Value = 8.589018950022027
$ws.Range("L4").Value = 12.26467696493445
$ws.Range("B5").Value = 19.19810666058401
$ws.Range("C5").Value = 9.331984054616353
$ws.Range("E5").Value = 14.40100566945615
$ws.Range("F5").Value = 41.49021884706757
$ws.Range("G5").Value = 43.31857704508758
$ws.Range("H5").Value = 17.85808548774006
$ws.Range("J5").Value = 8.596281143832018
$ws.Range("L5").Value = 12.2593204211205
$ws.Range("B6").Value = 19.17780667622644
$ws.Range("C6").Value = 9.318992067452998
$ws.Range("E6").Value = 14.4019317230173
$ws.Range("F6").Value = 41.49219350076076
$ws.Range("G6").Value = 43.31924510415857
$ws.Range("H6").Value = 17.86096521619932
$ws.Range("J6").Value = 8.597498678761978
$ws.Range("L6").Value = 12.25845445920616
$ws.Range("B7").Value = 19.31861872294413
$ws.Range("C7").Value = 9.408644158629281
$ws.Range("E7").Value = 14.39561330657002
$ws.Range("F7").Value = 41.47953430160437
$ws.Range("G7").Value = 43.31617932224619
$ws.Range("H7").Value = 17.84136536212236
$ws.Range("J7").Value = 8.589116109809213
$ws.Range("L7").Value = 12.26460315231061
$ws.Range("B8").Value = 19.93642532277051
$ws.Range("C8").Value = 9.790741882753755
$ws.Range("E8").Value = 14.37043456444942
$ws.Range("F8").Value = 41.44963710892702
$ws.Range("G8").Value = 43.34165974459078
$ws.Range("H8").Value = 17.76452409253688
$ws.Range("J8").Value = 8.553814173333526
$ws.Range("L8").Value = 12.29512534902488
$ws.Range("B9").Value = 21.12350790385494
$ws.Range("C9").Value = 10.48730447177641
$ws.Range("E9").Value = 14.33073407532871
$ws.Range("F9").Value = 41.4847198133956
$ws.Range("G9").Value = 43.53263632497294
$ws.Range("H9").Value = 17.64881838042121
$ws.Range("J9").Value = 8.490603528981712
$ws.Range("L9").Value = 12.36691014012361
$ws.Range("B10").Value = 21.97119325000723
$ws.Range("C10").Value = 10.96376968828047
$ws.Range("E10").Value = 14.30742822275552
$ws.Range("F10").Value = 41.5678193016087
$ws.Range("G10").Value = 43.75916525347935
$ws.Range("H10").Value = 17.5853990621897
$ws.Range("J10").Value = 8.4477801303261
$ws.Range("L10").Value = 12.42668354882264
$ws.Range("B11").Value = 22.34969097102511
$ws.Range("C11").Value = 11.172304923236
$ws.Range("E11").Value = 14.29809094554957
$ws.Range("F11").Value = 41.61806911086443
$ws.Range("G11").Value = 43.88087623935489
$ws.Range("H11").Value = 17.56128930507191
$ws.Range("J11").Value = 8.429073373950247
$ws.Range("L11").Value = 12.45534742876998
$ws.Range("B12").Value = 22.49185697068112
$ws.Range("C12").Value = 11.25005121939262
$ws.Range("E12").Value = 14.29473636896379
$ws.Range("F12").Value = 41.6388834238884
$ws.Range("G12").Value = 43.9296339801487
$ws.Range("H12").Value = 17.55284525009214
$ws.Range("J12").Value = 8.422100066046681
$ws.Range("L12").Value = 12.46640834573008
$ws.Range("B13").Value = 22.46129279908073
$ws.Range("C13").Value = 11.23336208750144
$ws.Range("E13").Value = 14.29545078582808
$ws.Range("F13").Value = 41.63432135430851
$ws.Range("G13").Value = 43.91901474995295
$ws.Range("H13").Value = 17.55463326354338
$ws.Range("J13").Value = 8.423596988363313
$ws.Range("L13").Value = 12.46401707958808
$ws.Range("B14").Value = 22.3614111178346
$ws.Range("C14").Value = 11.17872585453057
$ws.Range("E14").Value = 14.29781133356721
$ws.Range("F14").Value = 41.61974575828896
$ws.Range("G14").Value = 43.88483419542307
$ws.Range("H14").Value = 17.56058084021267
$ws.Range("J14").Value = 8.4284974646629
$ws.Range("L14").Value = 12.45625330868088
$ws.Range("B15").Value = 22.30007531759146
$ws.Range("C15").Value = 11.14509938417423
$ws.Range("E15").Value = 14.29928082275383
$ws.Range("F15").Value = 41.61105018875713
$ws.Range("G15").Value = 43.86424458194688
$ws.Range("H15").Value = 17.56431333277711
$ws.Range("J15").Value = 8.431513520856246
$ws.Range("L15").Value = 12.4515245118281
$ws.Range("B16").Value = 21.94630213908506
$ws.Range("C16").Value = 10.94997260214064
$ws.Range("E16").Value = 14.30806383219653
$ws.Range("F16").Value = 41.56478553036126
$ws.Range("G16").Value = 43.75158551579336
$ws.Range("H16").Value = 17.58707043462897
$ws.Range("J16").Value = 8.44901816882494
$ws.Range("L16").Value = 12.42483948743832
$ws.Range("B17").Value = 21.72734799007139
$ws.Range("C17").Value = 10.82813508254978
$ws.Range("E17").Value = 14.31377540486571
$ws.Range("F17").Value = 41.53958951302161
$ws.Range("G17").Value = 43.68724362582454
$ws.Range("H17").Value = 17.60224824873983
$ws.Range("J17").Value = 8.459954370585686
$ws.Range("L17").Value = 12.40884261150236
$ws.Range("B18").Value = 21.60074895222347
$ws.Range("C18").Value = 10.75728660989987
$ws.Range("E18").Value = 14.31717962435292
$ws.Range("F18").Value = 41.52626942925781
$ws.Range("G18").Value = 43.65199368394654
$ws.Range("H18").Value = 17.6114240411002
$ws.Range("J18").Value = 8.466317464530178
$ws.Range("L18").Value = 12.39978057203806
$ws.Range("B19").Value = 21.55777534583965
$ws.Range("C19").Value = 10.7331674749486
$ws.Range("E19").Value = 14.31835270411683
$ws.Range("F19").Value = 41.52196085913195
$ws.Range("G19").Value = 43.6403609322473
$ws.Range("H19").Value = 17.61460726227231
$ws.Range("J19").Value = 8.468484438578775
$ws.Range("L19").Value = 12.39673634051703
$ws.Range("B20").Value = 21.75072563235409
$ws.Range("C20").Value = 10.8411849743318
$ws.Range("E20").Value = 14.31315507889511
$ws.Range("F20").Value = 41.54215039781374
$ws.Range("G20").Value = 43.69391109492905
$ws.Range("H20").Value = 17.60058636322892
$ws.Range("J20").Value = 8.458782655849891
$ws.Range("L20").Value = 12.41053116373054
$ws.Range("B21").Value = 22.39078140437362
$ws.Range("C21").Value = 11.19480729227196
$ws.Range("E21").Value = 14.29711306937401
$ws.Range("F21").Value = 41.62397854183885
$ws.Range("G21").Value = 43.89480158092066
$ws.Range("H21").Value = 17.55881524845011
$ws.Range("J21").Value = 8.427055081988165
$ws.Range("L21").Value = 12.45852815406213
$ws.Range("B22").Value = 22.80225982451904
$ws.Range("C22").Value = 11.41878458564978
$ws.Range("E22").Value = 14.28768484979392
$ws.Range("F22").Value = 41.68786350973548
$ws.Range("G22").Value = 44.04163780703306
$ws.Range("H22").Value = 17.53551448843161
$ws.Range("J22").Value = 8.406963217681639
$ws.Range("L22").Value = 12.49109819021971
$ws.Range("B23").Value = 22.58331463902079
$ws.Range("C23").Value = 11.29990883237233
$ws.Range("E23").Value = 14.29262043302444
$ws.Range("F23").Value = 41.65281672395693
$ws.Range("G23").Value = 43.96185284483978
$ws.Range("H23").Value = 17.54758330075141
$ws.Range("J23").Value = 8.417627947180879
$ws.Range("L23").Value = 12.47360681985313
$ws.Range("B24").Value = 21.74015883028329
$ws.Range("C24").Value = 10.83528761466343
$ws.Range("E24").Value = 14.31343515270628
$ws.Range("F24").Value = 41.54098899126991
$ws.Range("G24").Value = 43.69089130781813
$ws.Range("H24").Value = 17.60133630067344
$ws.Range("J24").Value = 8.459312152135858
$ws.Range("L24").Value = 12.40976734907143
$ws.Range("B25").Value = 20.80601250569521
$ws.Range("C25").Value = 10.30485840935989
$ws.Range("E25").Value = 14.34044203332389
$ws.Range("F25").Value = 41.46516810436637
$ws.Range("G25").Value = 43.46582733582503
$ws.Range("H25").Value = 17.67634697630739
$ws.Range("J25").Value = 8.507064797010335
$ws.Range("L25").Value = 12.34623704186154
